$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item("Google Shape;122;p1")
$tr = $sh.TextFrame.TextRange

# Append a new bullet paragraph after the last existing one, preserving the
# inherited run formatting (Calibri / #00B0F0) that PowerPoint carries over
# from the preceding "TeleOp" run. Insert it as three separate runs so the
# word "Tophat" stays its own run, matching how "Auton"/"TeleOp" are split
# elsewhere in this shape.
$null = $tr.InsertAfter("`rLeft Trigger, Right Trigger & A to put robot ")
$null = $tr.InsertAfter("Tophat")
$null = $tr.InsertAfter(" in navigation mode")

# The text box auto-fits to its content (a:spAutoFit); growing by one bullet
# line increases its height. Match the resulting box size explicitly.
$sh.Height = 196.2953
